$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 26; this pushes the former rows 26-33 down
# to 27-34, and grows the sheet dimension from A1:R33 to A1:R34.
$ws.Rows.Item(26).Insert()

# Populate the newly inserted row 26 with the new weekly price entry.
$ws.Cells.Item(26,1).Value  = 5
$ws.Cells.Item(26,2).Value  = "Macroferia Regional de Talca"
$ws.Cells.Item(26,3).Value  = "Maule"
$ws.Cells.Item(26,4).Value  = 44736
$ws.Cells.Item(26,5).Value  = 7
$ws.Cells.Item(26,6).Value  = 100112043
$ws.Cells.Item(26,7).Value  = "Pepino dulce"
$ws.Cells.Item(26,8).Value  = "Cultivar IV Región"
$ws.Cells.Item(26,9).Value  = "Primera"
$ws.Cells.Item(26,10).Value = 200
$ws.Cells.Item(26,11).Value = 15000
$ws.Cells.Item(26,12).Value = 15000
$ws.Cells.Item(26,13).Value = 15000
$ws.Cells.Item(26,14).Value = '$/bandeja 18 kilos'
$ws.Cells.Item(26,15).Value = "Provincia de Limarí"
$ws.Cells.Item(26,16).Value = 833
$ws.Cells.Item(26,17).Value = 18
$ws.Cells.Item(26,18).Value = "Hortaliza"
